$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "Not equals to 0"
$ws.Range("A1").Value = "Conditional formatting on row (red background if condition true):"
$ws.Range("A3").Value = "Equals to the cell above"
$ws.Range("A4").Value = "NOT equals to the cell above"
$ws.Range("A6").Value = "Equals to zero"
$ws.Range("A7").Value = "Not equals to boolean FALSE"

$ws.Range("A3:A4").Font.Color = 0

$ws.Columns.Item(1).ColumnWidth = 57.83

$ws.Range("B9").Select() | Out-Null
